# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets:
#   F2: 90 -> 96
#   F3: 14 -> 15

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 96
    $ws.Range("F3").Value = 15
}
